# Update "Förändrad" (changed) date column C for existing rows 2..82
# from 45184 (2023-09-15) to 45186 (2023-09-17), and append two new
# data rows (83 and 84) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 82.
for ($r = 2; $r -le 82; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45186
    }
}

$dateFormat = $ws.Cells.Item(2, 2).NumberFormat

# Row 83 - new record
$ws.Cells.Item(83, 1).Value = "A 43549-2023"
$ws.Cells.Item(83, 2).Value2 = 45184
$ws.Cells.Item(83, 3).Value2 = 45186
$ws.Cells.Item(83, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item(83, 5).Value = "MUNKFORS"
$ws.Cells.Item(83, 6).Value = "Bergvik skog väst AB"
$ws.Cells.Item(83, 7).Value = 2.6
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = 0
$ws.Cells.Item(83, 14).Value = 0
$ws.Cells.Item(83, 15).Value = 0
$ws.Cells.Item(83, 16).Value = 0
$ws.Cells.Item(83, 17).Value = 0
$ws.Cells.Item(83, 2).NumberFormat = $dateFormat
$ws.Cells.Item(83, 3).NumberFormat = $dateFormat
$ws.Cells.Item(83, 18).WrapText = $true

# Row 84 - new record
$ws.Cells.Item(84, 1).Value = "A 43554-2023"
$ws.Cells.Item(84, 2).Value2 = 45184
$ws.Cells.Item(84, 3).Value2 = 45186
$ws.Cells.Item(84, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item(84, 5).Value = "MUNKFORS"
$ws.Cells.Item(84, 6).Value = "Bergvik skog väst AB"
$ws.Cells.Item(84, 7).Value = 2.5
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(84, 14).Value = 0
$ws.Cells.Item(84, 15).Value = 0
$ws.Cells.Item(84, 16).Value = 0
$ws.Cells.Item(84, 17).Value = 0
$ws.Cells.Item(84, 2).NumberFormat = $dateFormat
$ws.Cells.Item(84, 3).NumberFormat = $dateFormat
$ws.Cells.Item(84, 18).WrapText = $true

# Row heights: rows 82 and 83 have explicit 15pt custom height set,
# row 84 keeps the default (no explicit custom height).
$ws.Rows.Item(82).RowHeight = 15
$ws.Rows.Item(83).RowHeight = 15
